$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Network"
$ws.Range("D7").Value = "Server->Client"
$ws.Range("D8").Value = "Client->Server"
